# Restored from revision of admin on 01/01/2021 09:03:40 AM.TEST Author: admin. Type: SAVE.
# The only substantive content change in this revision is the value of
# cell C10 on the "Rules" sheet, which goes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
